$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-OrEmpty($val) {
    if ($val -eq $null) { return "" }
    return $val
}

# --- Rows 16, 17, 18: rotate data (new16=old18, new17=old16, new18=old17) ---
$old16 = @{ A = $ws.Range("A16").Value2; Q = $ws.Range("Q16").Value2; R = $ws.Range("R16").Value2; S = $ws.Range("S16").Value2; Z = $ws.Range("Z16").Value2; AB = $ws.Range("AB16").Value2; AC = Get-OrEmpty $ws.Range("AC16").Value2 }
$old17 = @{ A = $ws.Range("A17").Value2; Q = $ws.Range("Q17").Value2; R = $ws.Range("R17").Value2; S = $ws.Range("S17").Value2; Z = $ws.Range("Z17").Value2; AB = $ws.Range("AB17").Value2; AC = Get-OrEmpty $ws.Range("AC17").Value2 }
$old18 = @{ A = $ws.Range("A18").Value2; Q = $ws.Range("Q18").Value2; R = $ws.Range("R18").Value2; S = $ws.Range("S18").Value2; Z = $ws.Range("Z18").Value2; AB = $ws.Range("AB18").Value2; AC = Get-OrEmpty $ws.Range("AC18").Value2 }

$ws.Range("A16").Value = $old18.A
$ws.Range("Q16").Value = $old18.Q
$ws.Range("R16").Value = $old18.R
$ws.Range("S16").Value = $old18.S
$ws.Range("Z16").Value = $old18.Z
$ws.Range("AB16").Value = $old18.AB
$ws.Range("AC16").Value = $old18.AC

$ws.Range("A17").Value = $old16.A
$ws.Range("Q17").Value = $old16.Q
$ws.Range("R17").Value = $old16.R
$ws.Range("S17").Value = $old16.S
$ws.Range("Z17").Value = $old16.Z
$ws.Range("AB17").Value = $old16.AB
$ws.Range("AC17").Value = $old16.AC

$ws.Range("A18").Value = $old17.A
$ws.Range("Q18").Value = $old17.Q
$ws.Range("R18").Value = $old17.R
$ws.Range("S18").Value = $old17.S
$ws.Range("Z18").Value = $old17.Z
$ws.Range("AB18").Value = $old17.AB
$ws.Range("AC18").Value = $old17.AC

# --- Rows 20, 21: swap data ---
$old20 = @{ A = $ws.Range("A20").Value2; Q = $ws.Range("Q20").Value2; R = $ws.Range("R20").Value2; S = $ws.Range("S20").Value2; Z = $ws.Range("Z20").Value2; AB = $ws.Range("AB20").Value2 }
$old21 = @{ A = $ws.Range("A21").Value2; Q = $ws.Range("Q21").Value2; R = $ws.Range("R21").Value2; S = $ws.Range("S21").Value2; Z = $ws.Range("Z21").Value2; AB = $ws.Range("AB21").Value2 }

$ws.Range("A20").Value = $old21.A
$ws.Range("Q20").Value = $old21.Q
$ws.Range("R20").Value = $old21.R
$ws.Range("S20").Value = $old21.S
$ws.Range("Z20").Value = $old21.Z
$ws.Range("AB20").Value = $old21.AB

$ws.Range("A21").Value = $old20.A
$ws.Range("Q21").Value = $old20.Q
$ws.Range("R21").Value = $old20.R
$ws.Range("S21").Value = $old20.S
$ws.Range("Z21").Value = $old20.Z
$ws.Range("AB21").Value = $old20.AB

# --- Rows 22, 23: swap data ---
$old22 = @{ A = $ws.Range("A22").Value2; B = $ws.Range("B22").Value2; E = $ws.Range("E22").Value2; F = $ws.Range("F22").Value2; G = $ws.Range("G22").Value2; H = $ws.Range("H22").Value2; Q = $ws.Range("Q22").Value2; R = $ws.Range("R22").Value2; S = $ws.Range("S22").Value2; Z = $ws.Range("Z22").Value2; AB = $ws.Range("AB22").Value2 }
$old23 = @{ A = $ws.Range("A23").Value2; B = $ws.Range("B23").Value2; E = $ws.Range("E23").Value2; F = $ws.Range("F23").Value2; G = $ws.Range("G23").Value2; H = $ws.Range("H23").Value2; Q = $ws.Range("Q23").Value2; R = $ws.Range("R23").Value2; S = $ws.Range("S23").Value2; Z = $ws.Range("Z23").Value2; AB = $ws.Range("AB23").Value2 }

$ws.Range("A22").Value = $old23.A
$ws.Range("B22").Value = $old23.B
$ws.Range("E22").Value = $old23.E
$ws.Range("F22").Value = $old23.F
$ws.Range("G22").Value = $old23.G
$ws.Range("H22").Value = $old23.H
$ws.Range("Q22").Value = $old23.Q
$ws.Range("R22").Value = $old23.R
$ws.Range("S22").Value = $old23.S
$ws.Range("Z22").Value = $old23.Z
$ws.Range("AB22").Value = $old23.AB

$ws.Range("A23").Value = $old22.A
$ws.Range("B23").Value = $old22.B
$ws.Range("E23").Value = $old22.E
$ws.Range("F23").Value = $old22.F
$ws.Range("G23").Value = $old22.G
$ws.Range("H23").Value = $old22.H
$ws.Range("Q23").Value = $old22.Q
$ws.Range("R23").Value = $old22.R
$ws.Range("S23").Value = $old22.S
$ws.Range("Z23").Value = $old22.Z
$ws.Range("AB23").Value = $old22.AB

# --- Rows 25, 27: swap data ---
$old25 = @{ A = $ws.Range("A25").Value2; B = $ws.Range("B25").Value2; E = $ws.Range("E25").Value2; F = $ws.Range("F25").Value2; G = $ws.Range("G25").Value2; H = $ws.Range("H25").Value2; Q = $ws.Range("Q25").Value2; R = $ws.Range("R25").Value2; S = $ws.Range("S25").Value2; Z = $ws.Range("Z25").Value2; AB = $ws.Range("AB25").Value2 }
$old27 = @{ A = $ws.Range("A27").Value2; B = $ws.Range("B27").Value2; E = $ws.Range("E27").Value2; F = $ws.Range("F27").Value2; G = $ws.Range("G27").Value2; H = $ws.Range("H27").Value2; Q = $ws.Range("Q27").Value2; R = $ws.Range("R27").Value2; S = $ws.Range("S27").Value2; Z = $ws.Range("Z27").Value2; AB = $ws.Range("AB27").Value2 }

$ws.Range("A25").Value = $old27.A
$ws.Range("B25").Value = $old27.B
$ws.Range("E25").Value = $old27.E
$ws.Range("F25").Value = $old27.F
$ws.Range("G25").Value = $old27.G
$ws.Range("H25").Value = $old27.H
$ws.Range("Q25").Value = $old27.Q
$ws.Range("R25").Value = $old27.R
$ws.Range("S25").Value = $old27.S
$ws.Range("Z25").Value = $old27.Z
$ws.Range("AB25").Value = $old27.AB

$ws.Range("A27").Value = $old25.A
$ws.Range("B27").Value = $old25.B
$ws.Range("E27").Value = $old25.E
$ws.Range("F27").Value = $old25.F
$ws.Range("G27").Value = $old25.G
$ws.Range("H27").Value = $old25.H
$ws.Range("Q27").Value = $old25.Q
$ws.Range("R27").Value = $old25.R
$ws.Range("S27").Value = $old25.S
$ws.Range("Z27").Value = $old25.Z
$ws.Range("AB27").Value = $old25.AB
